$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume/1h (E) updates
$deData = @{
    2 = @("303.88", "0.10%")
    3 = @("35.86", "-3.26%")
    4 = @("5.063", "1.13%")
    5 = @("0.07888", "0.34%")
    6 = @("2.109", "-4.55%")
    7 = @("7.944", "-0.67%")
    8 = @("0.9213", "0.20%")
    9 = @("0.09712", "1.43%")
    10 = @("0.1836", "-2.59%")
    11 = @("0.08607", "-0.01%")
    12 = @("0.03579", "-0.63%")
    13 = @("0.09934", "-0.47%")
    14 = @("0.001435", "-2.97%")
    15 = @("0.005745", "2.01%")
    16 = @("3.464", "0.37%")
    17 = @($null, "2.81%")
    18 = @("2.751", "22.31%")
    19 = @("0.3380", "-1.15%")
    20 = @("0.1349", "2.55%")
    21 = @("5.183", "8.95%")
    22 = @("0.2216", "0.76%")
    23 = @("0.04549", "-1.05%")
    24 = @("0.001238", "0.64%")
    25 = @("0.004866", "8.97%")
    26 = @("0.0001305", "-6.74%")
    27 = @("0.0004767", "0.34%")
    39 = @("0.01837", "-0.65%")
    40 = @("0.04701", "-0.90%")
    41 = @("0.007909", "-2.64%")
    42 = @("0.1392", "-0.14%")
    43 = @("0.007595", "0.58%")
    44 = @("0.002198", "-0.52%")
    45 = @("0.01125", "6.92%")
    46 = @("0.00006303", "-1.13%")
    47 = @("0.00000000753", "0.40%")
    48 = @($null, "0.23%")
    49 = @("50.54", "89.28%")
    50 = @("0.001907", "-29.11%")
    51 = @("0.00002108", "0.40%")
}

foreach ($row in $deData.Keys) {
    $vals = $deData[$row]
    if ($vals[0] -ne $null) {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $vals[0]
    }
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Hora (G) column: 13 -> 14 for all data rows 2..51
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = "14"
}
